$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capitalise the team-member names in the header row (C2:H2)
$ws.Range("C2").Value = "Ravi"
$ws.Range("D2").Value = "Ben"
$ws.Range("E2").Value = "Zach"
$ws.Range("F2").Value = "Harry"
$ws.Range("G2").Value = "Freddie"
$ws.Range("H2").Value = "Adam"

# 2) Add a new meeting row (row 6: 11 Feb 2020)
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 43872

$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = "15:00 - 16:10"

$ws.Range("C3").Copy()
$ws.Range("C6:G6").PasteSpecial(-4122)

$ws.Range("H6").Font.Color = 24832

$ws.Range("J6").Value = "Adam - family situation"

# 3) Note for the earlier meeting (row 5) about Zach's absence
$ws.Range("J5").Value = "Zach - assessment centre"

# 4) Rename the two header titles
$ws.Range("C1").Value = "Team Members"
$ws.Range("B2").Value = "Times"

$ws.Application.CutCopyMode = $false

# 5) Update the view state (scrolled right one column, selection on G16)
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
[void]$ws.Range("G16").Select()
